# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price column values are purely numeric-looking text (e.g. "22.80") that
# Excel would otherwise auto-convert to a number (dropping the trailing zero).
# Pre-format those specific cells as Text so the assigned values stay as strings.
$textCells = @("D5", "D8", "D18", "D19", "D25", "D26", "D30", "D36", "D37", "D42", "D46", "D47", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.138.00'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.678.39'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '214.22'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '22.80'
$ws.Range('E8').Value = '  +7.01%  '
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.916.45'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = '1.676.11'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  +2.40%  '
$ws.Range('E15').Value = '  +3.39%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '27.103.16'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '235.69'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').Value = '7.90'
$ws.Range('E19').Value = '  -2.58%  '
$ws.Range('D20').Value = '0.0₃0741'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('E23').Value = '  +2.87%  '
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').Value = '147.21'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = '7.42'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '0.0501'
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').Value = '1.541.89'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('E35').Value = '  -2.90%  '
$ws.Range('D36').Value = '0.607'
$ws.Range('E36').Value = '  +3.34%  '
$ws.Range('D37').Value = '0.940'
$ws.Range('E37').Value = '  +2.41%  '
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  +2.69%  '
$ws.Range('E41').Value = '  +3.24%  '
$ws.Range('D42').Value = '69.59'
$ws.Range('E42').Value = '  +2.45%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '1.822.67'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').Value = '0.779'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('D47').Value = '89.92'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('D49').Value = '1.63'
$ws.Range('E49').Value = '  +6.14%  '
$ws.Range('D50').Value = '8.21'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('E51').Value = '  -0.10%  '
